# Refresh the "Elapsed Duration(Hrs)" column (G) on several of the R1..R6
# outage-report sheets to reflect a later snapshot time.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3918:35:23"
$ws1.Range("G3").Value = "58:08:01"

$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12099:59:04"
$ws2.Range("G3").Value = "3229:42:33"
$ws2.Range("G4").Value = "467:54:07"

$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2945:48:53"
$ws4.Range("G3").Value = "173:01:08"

$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "419:47:52"

$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "60:20:10"
